$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# The trailing footnote currently sitting in B98 needs to move down to B99
# so the new data row can be inserted at row 98.
$ws.Range("B99").Value = $ws.Range("B98").Value2

# Copy the formatting of the last data row (97) down into the new row (98)
# so the new cells pick up the same number formats/styles as the rest of
# the table, then fill in the new day's figures.
$ws.Range("A97:E97").Copy()
$ws.Range("A98:E98").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A98").Value = "2020-05-02"
$ws.Range("B98").Value = 292
$ws.Range("C98").Value = 32728
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 6958

# Extend the print area to include the newly-added row (table now ends at
# row 100 once the moved footnote row is accounted for).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = '=相談件数!$A$1:$E$100'
    }
}

# Reflect the new active selection on the sheet.
$ws.Range("C100").Select()
